$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.154.33"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.484.01"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.88"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.63"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.42%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.68%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +5.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.124"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.388"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.072.64"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.58%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.04%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.94%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.482.08"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.116.56"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.48"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -6.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.01"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.54"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "385.11"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.620.60"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.41"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.37%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.25%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.19%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -4.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.45"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.96"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.94%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.509.57"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.40%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.04"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.22"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.82"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.13%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "163.57"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.69%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0781"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.42%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.53%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.35"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "23.97"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -6.75%  "
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.63"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.36%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.15"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.91%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.940"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +5.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.76"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.368.82"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.75%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.31%  "
